$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D7").Value = -7.820099999999999
$ws.Range("A9").Value = -21.92359999999998
$ws.Range("D12").Value = -7.123699999999998
$ws.Range("D14").Value = -7.654100000000003
$ws.Range("A18").Value = -22.26060000000001
$ws.Range("A20").Value = -19.26499999999998
$ws.Range("D26").Value = -8.533099999999997
$ws.Range("A27").Value = -22.01789999999999
$ws.Range("D27").Value = -8.786599999999998
$ws.Range("D29").Value = -7.331600000000001
$ws.Range("A35").Value = -19.13519999999999
$ws.Range("D37").Value = -7.576800000000002
$ws.Range("D38").Value = -8.190499999999995
$ws.Range("D51").Value = -7.642499999999999
$ws.Range("D52").Value = -7.6674
$ws.Range("D55").Value = -8.973399999999996
$ws.Range("A69").Value = -21.5904
$ws.Range("D69").Value = -7.201099999999996
$ws.Range("D70").Value = -7.345299999999998
$ws.Range("A76").Value = -19.74259999999998
$ws.Range("A78").Value = -20.17219999999998
$ws.Range("D81").Value = -7.544900000000001
$ws.Range("A82").Value = -21.78530000000001
$ws.Range("A83").Value = -22.1355
$ws.Range("D83").Value = -8.380500000000005
$ws.Range("A93").Value = -21.20869999999999
$ws.Range("D102").Value = -7.607700000000002
